# Adds a new "2022-Q4" sheet (with fund holding data) positioned between
# "总计" and "2021-Q2", and records a summary row for it on "总计".

function Set-ExactText {
    # Writes $text into $range as a literal text value, even when it looks
    # like a number (e.g. "009225", "94.90"), without leaving behind any
    # extra/applied number-format style on the cell.
    param($range, [string]$text)
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)   # "总计"
$q2Sheet    = $wb.Worksheets.Item(2)   # "2021-Q2"

# Insert the new sheet right before "2021-Q2" (so order becomes 总计, 2022-Q4, 2021-Q2)
$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q4"

# ---- Build "2022-Q4" sheet ----

# Copy the bold/bordered header style (used on 总计!B1:D1) onto the new header row
$totalSheet.Range("B1:D1").Copy()
$newSheet.Range("B1:D1").PasteSpecial(-4122)   # xlPasteFormats
$newSheet.Range("E1:H1").PasteSpecial(-4122)   # xlPasteFormats

# Copy the bold/bordered index-column style (used on 总计!A2) onto A2:A6
$totalSheet.Range("A2").Copy()
$newSheet.Range("A2:A6").PasteSpecial(-4122)   # xlPasteFormats

# Header labels
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Index column (A) and rank column (H) are genuine numbers
$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1
$newSheet.Range("A4").Value = 2
$newSheet.Range("A5").Value = 3
$newSheet.Range("A6").Value = 4

$newSheet.Range("H2").Value = 5
$newSheet.Range("H3").Value = 3
$newSheet.Range("H4").Value = 3
$newSheet.Range("H5").Value = 3
$newSheet.Range("H6").Value = 5

# Fund code / name / figures are stored as text (preserves leading/trailing zeros)
Set-ExactText $newSheet.Range("B2") "009225"
Set-ExactText $newSheet.Range("C2") "天弘中证中美互联网指数（QDII）A"
Set-ExactText $newSheet.Range("D2") "1.27"
Set-ExactText $newSheet.Range("E2") "94.90"
Set-ExactText $newSheet.Range("F2") "9.16"
Set-ExactText $newSheet.Range("G2") "0.1163"

Set-ExactText $newSheet.Range("B3") "012751"
Set-ExactText $newSheet.Range("C3") "建信纳斯达克100指数（QDII）A 美元现汇"
Set-ExactText $newSheet.Range("D3") "1.06"
Set-ExactText $newSheet.Range("E3") "82.28"
Set-ExactText $newSheet.Range("F3") "5.59"
Set-ExactText $newSheet.Range("G3") "0.0593"

Set-ExactText $newSheet.Range("B4") "012752"
Set-ExactText $newSheet.Range("C4") "建信纳斯达克100指数（QDII）C 人民币"
Set-ExactText $newSheet.Range("D4") "1.06"
Set-ExactText $newSheet.Range("E4") "82.28"
Set-ExactText $newSheet.Range("F4") "5.59"
Set-ExactText $newSheet.Range("G4") "0.0593"

Set-ExactText $newSheet.Range("B5") "012753"
Set-ExactText $newSheet.Range("C5") "建信纳斯达克100指数（QDII）C 美元现汇"
Set-ExactText $newSheet.Range("D5") "1.06"
Set-ExactText $newSheet.Range("E5") "82.28"
Set-ExactText $newSheet.Range("F5") "5.59"
Set-ExactText $newSheet.Range("G5") "0.0593"

Set-ExactText $newSheet.Range("B6") "009226"
Set-ExactText $newSheet.Range("C6") "天弘中证中美互联网指数（QDII）C"
Set-ExactText $newSheet.Range("D6") "0.63"
Set-ExactText $newSheet.Range("E6") "94.90"
Set-ExactText $newSheet.Range("F6") "9.16"
Set-ExactText $newSheet.Range("G6") "0.0577"

# ---- Update "总计" sheet: add a 2022-Q4 row, push 2021-Q2 row down ----

# Duplicate the existing data row's formatting into row 3 before overwriting row 2
$totalSheet.Range("A2:D2").Copy()
$totalSheet.Range("A3:D3").PasteSpecial(-4122)  # xlPasteFormats

# Row 3 keeps the former "2021-Q2" figures
$totalSheet.Range("A3").Value = 1
Set-ExactText $totalSheet.Range("B3") "2021-Q2"
$totalSheet.Range("C3").Value = 7
$totalSheet.Range("D3").Value = 0.9399999999999999

# Row 2 now holds the new "2022-Q4" figures
$totalSheet.Range("A2").Value = 0
Set-ExactText $totalSheet.Range("B2") "2022-Q4"
$totalSheet.Range("C2").Value = 5
$totalSheet.Range("D2").Value = 0.35
